$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts the existing
# Jun_17 / Jun_15 / Jun_13 / Jun_10 columns (B:E) one place to the
# right (to C:F) and leaves a fresh, empty column B in their place.
$ws.Columns("B").Insert()

# New "Jun_26" week column header.
$ws.Range("B1").Value = "Jun_26"

# Fill the new column with the "UN" rating used throughout the sheet.
$ws.Range("B2:B27").Value = "UN"

# Append two new tickers at the bottom of the table.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
